$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 33.0671278619396
$ws.Range("B4").Value = 0.009999990463256836
$ws.Range("B6").Value = 33.0671278619396
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 13
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 12
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 9
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 7

$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 32.61192465059682
$ws.Range("B4").Value = 34.69770569366315
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 30
$ws.Range("B7").Value = 36.71579249669672
$ws.Range("B8").Value = 30.34885527085025
$ws.Range("B9").Value = 32.01159140980468
$ws.Range("B10").Value = 32.31224998648503
$ws.Range("B11").Value = 34.76592070603971
$ws.Range("B12").Value = 30
$ws.Range("B13").Value = 37.94859027624736
$ws.Range("B14").Value = 37.90090852477161
$ws.Range("B15").Value = 37.27819014430416

$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 250.970000000001
$ws.Range("C8").Value = 260.990000000001
$ws.Range("C9").Value = 252.9750000000009
$ws.Range("C10").Value = 269.580000000001
$ws.Range("C11").Value = 250.575000000001
$ws.Range("C12").Value = 235.775
$ws.Range("C13").Value = 229.025
$ws.Range("C14").Value = 213.42
$ws.Range("C15").Value = 226.76
$ws.Range("C16").Value = 221.56
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C22").Value = 72.6299999999995
$ws.Range("C23").Value = 80.0549999999995
$ws.Range("C24").Value = 82.31999999999948
$ws.Range("C25").Value = 83.9549999999995
$ws.Range("C26").Value = 80.8149999999995
$ws.Range("C27").Value = 295.9199999999997
$ws.Range("C28").Value = 323.5
$ws.Range("C29").Value = 294.2649999999996
$ws.Range("C30").Value = 311.1
$ws.Range("C31").Value = 297.3649999999997
$ws.Range("C32").Value = 154.3
$ws.Range("C33").Value = 148.3449999999993
$ws.Range("C34").Value = 128.7049999999993
$ws.Range("C35").Value = 146.3249999999993
$ws.Range("C36").Value = 134.2149999999993
$ws.Range("C37").Value = 141.0250000000001
$ws.Range("C38").Value = 143.4
$ws.Range("C39").Value = 139.7050000000001
$ws.Range("C40").Value = 150.4250000000002
$ws.Range("C41").Value = 134.7700000000002
$ws.Range("C42").Value = 140.5549999999989
$ws.Range("C43").Value = 159.2149999999988
$ws.Range("C44").Value = 142.1399999999988
$ws.Range("C45").Value = 147.7249999999989
$ws.Range("C46").Value = 139.7449999999989
$ws.Range("C47").Value = 226.0399999999994
$ws.Range("C48").Value = 247.1799999999993
$ws.Range("C49").Value = 221.8549999999994
$ws.Range("C50").Value = 238.4549999999994
$ws.Range("C51").Value = 224.4749999999994
$ws.Range("C52").Value = 57.95
$ws.Range("C53").Value = 58.67999999999927
$ws.Range("C54").Value = 61.72999999999927
$ws.Range("C55").Value = 60.65499999999928
$ws.Range("C56").Value = 52.91499999999927
$ws.Range("C57").Value = 250.970000000001
$ws.Range("C58").Value = 260.990000000001
$ws.Range("C59").Value = 252.9750000000009
$ws.Range("C60").Value = 269.580000000001
$ws.Range("C61").Value = 250.575000000001
$ws.Range("C62").Value = 295.9199999999997
$ws.Range("C67").Value = 235.775
$ws.Range("C68").Value = 229.025
$ws.Range("C69").Value = 213.42
$ws.Range("C70").Value = 226.76
$ws.Range("C71").Value = 221.56

$ws = $wb.Worksheets.Item("R")
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0

# Delete rows 2-4 on sheet "alpha" (dimension A1:C4 -> A1:C1)
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2:C4").EntireRow.Delete()

# Delete rows 2-4 on sheet "y" (dimension A1:D4 -> A1:D1)
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2:D4").EntireRow.Delete()

# Delete rows 2-6 on sheet "rho" (dimension A1:C6 -> A1:C1)
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2:C6").EntireRow.Delete()
